# Daily update at 8 AM UTC: append today's snapshot as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 30
$newRow = $lastRow + 1

# Carry forward today's totals (same date, same cumulative win counts as
# the prior snapshot captured earlier the same day).
$ws.Cells.Item($newRow, 1).Value = $ws.Cells.Item($lastRow, 1).Value2
$ws.Cells.Item($newRow, 2).Value = $ws.Cells.Item($lastRow, 2).Value2
$ws.Cells.Item($newRow, 3).Value = $ws.Cells.Item($lastRow, 3).Value2
$ws.Cells.Item($newRow, 4).Value = $ws.Cells.Item($lastRow, 4).Value2

# The "most recent row" date cell is styled as a date-only value (no time);
# older rows show full date+time. Move that styling down to the new last row,
# and restyle the previous last row like all the other historical rows.
$ws.Cells.Item($lastRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD"
